# Actualización automática 2025-07-04 17:10:08
#
# The advisor "HIDALGO HIDALGO PEDRO GUSTAVO" registers a new PORCELANATO
# sale of 3233.78 for client "JARAMILLO CARVAJAL NICOLAS ESTEBAN" in julio.
# This propagates into the per-group sheet, the monthly-sales sheet and the
# monthly-compliance summary sheet (plus their totals / ratios / counters).

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Sheet "VENTAS POR GRUPO": new PORCELANATO sale for row 10 (client
# JARAMILLO CARVAJAL NICOLAS ESTEBAN), and bump the "x de 20" counter
# for the PORCELANATO column (M) from 2 to 3 buyers.
# ------------------------------------------------------------------
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsGrupo.Cells.Item(10, 13).Value = 3233.78
$wsGrupo.Cells.Item(22, 13).Value = "3 de 20"

# ------------------------------------------------------------------
# Sheet "VENTA MENSUAL": the same sale lands in julio (column F) for
# the same client row, and the julio total (row 22) grows accordingly.
# ------------------------------------------------------------------
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsMensual.Cells.Item(10, 6).Value = 3233.78
$wsMensual.Cells.Item(22, 6).Value = 13495.38

# ------------------------------------------------------------------
# Sheet "CUMPLIMIENTO MENSUAL": PORCELANATO (row 16) VENTA / POR CUMPLIR /
# CUMPLIMIENTO values update, as does the TOTAL row (19). The VENTA
# column (D) also grows wide enough that Excel bumps its stored column
# width from 13 to 14 characters.
# ------------------------------------------------------------------
$wsCumpl = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

$wsCumpl.Columns.Item(4).ColumnWidth = 13.1667

$wsCumpl.Cells.Item(16, 4).Value = 10541.24
$wsCumpl.Cells.Item(16, 5).Value = 33725
$wsCumpl.Cells.Item(16, 6).Value = 0.23813271694185

$wsCumpl.Cells.Item(19, 4).Value = 13495.38
$wsCumpl.Cells.Item(19, 5).Value = 51882.61762291769
$wsCumpl.Cells.Item(19, 6).Value = 0.2064208218464818
